$d = $word.ActiveDocument

$d.Content.Find.Execute("320×6=1920", $true, $false, $false, $false, $false, $true, 1, $false, "717×2=1434", 2) | Out-Null
$d.Content.Find.Execute("234×2=468", $true, $false, $false, $false, $false, $true, 1, $false, "285×9=2565", 2) | Out-Null
$d.Content.Find.Execute("884×3=2652", $true, $false, $false, $false, $false, $true, 1, $false, "845×7=5915", 2) | Out-Null
$d.Content.Find.Execute("572×8=4576", $true, $false, $false, $false, $false, $true, 1, $false, "643×5=3215", 2) | Out-Null
$d.Content.Find.Execute("164×7=1148", $true, $false, $false, $false, $false, $true, 1, $false, "888×6=5328", 2) | Out-Null
$d.Content.Find.Execute("224×7=1568", $true, $false, $false, $false, $false, $true, 1, $false, "889×2=1778", 2) | Out-Null
$d.Content.Find.Execute("301×3=903", $true, $false, $false, $false, $false, $true, 1, $false, "410×3=1230", 2) | Out-Null
$d.Content.Find.Execute("512×7=3584", $true, $false, $false, $false, $false, $true, 1, $false, "962×7=6734", 2) | Out-Null
$d.Content.Find.Execute("221×8=1768", $true, $false, $false, $false, $false, $true, 1, $false, "734×5=3670", 2) | Out-Null
$d.Content.Find.Execute("348×7=2436", $true, $false, $false, $false, $false, $true, 1, $false, "589×2=1178", 2) | Out-Null
$d.Content.Find.Execute("847×5=4235", $true, $false, $false, $false, $false, $true, 1, $false, "132×8=1056", 2) | Out-Null
$d.Content.Find.Execute("812×4=3248", $true, $false, $false, $false, $false, $true, 1, $false, "450×6=2700", 2) | Out-Null
$d.Content.Find.Execute("509×8=4072", $true, $false, $false, $false, $false, $true, 1, $false, "305×2=610", 2) | Out-Null
$d.Content.Find.Execute("312×4=1248", $true, $false, $false, $false, $false, $true, 1, $false, "657×7=4599", 2) | Out-Null
$d.Content.Find.Execute("316×8=2528", $true, $false, $false, $false, $false, $true, 1, $false, "977×6=5862", 2) | Out-Null
$d.Content.Find.Execute("552×6=3312", $true, $false, $false, $false, $false, $true, 1, $false, "979×9=8811", 2) | Out-Null
$d.Content.Find.Execute("855×4=3420", $true, $false, $false, $false, $false, $true, 1, $false, "251×3=753", 2) | Out-Null
$d.Content.Find.Execute("507×4=2028", $true, $false, $false, $false, $false, $true, 1, $false, "884×6=5304", 2) | Out-Null
$d.Content.Find.Execute("869×9=7821", $true, $false, $false, $false, $false, $true, 1, $false, "202×6=1212", 2) | Out-Null
$d.Content.Find.Execute("640×4=2560", $true, $false, $false, $false, $false, $true, 1, $false, "275×5=1375", 2) | Out-Null
$d.Content.Find.Execute("371×2=742", $true, $false, $false, $false, $false, $true, 1, $false, "466×3=1398", 2) | Out-Null
$d.Content.Find.Execute("307×9=2763", $true, $false, $false, $false, $false, $true, 1, $false, "839×9=7551", 2) | Out-Null
$d.Content.Find.Execute("264×7=1848", $true, $false, $false, $false, $false, $true, 1, $false, "923×6=5538", 2) | Out-Null
$d.Content.Find.Execute("842×7=5894", $true, $false, $false, $false, $false, $true, 1, $false, "247×7=1729", 2) | Out-Null
$d.Content.Find.Execute("578×8=4624", $true, $false, $false, $false, $false, $true, 1, $false, "831×7=5817", 2) | Out-Null
